$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a third row of signup data (phone number + hyperlinked name) ---
$ws.Range("A3").Value = 9843801062

# Create the hyperlink (mailto, same convention as the existing
# B2 -> Velskar@1 link) with its own display text ("vedha@1") ...
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:vedha@1", [Type]::Missing, [Type]::Missing, "vedha@1")

# ... then set the cell's visible text to the actual name, and re-apply the
# built-in "Hyperlink" style so B3 matches B2's formatting.
$ws.Range("B3").Value = "karthiga"
$ws.Range("B3").Style = "Hyperlink"

# Reflect the new active cell in the sheet's selection.
$ws.Range("B3").Select()
